$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.738607
$ws.Range("H2").Value = 8.215821
$ws.Range("I2").Value = 0.2235648590725649
$ws.Range("J2").Value = 0.223564859072565
$ws.Range("M2").Value = 0.5942546666666666
$ws.Range("N2").Value = 1.782764
$ws.Range("O2").Value = 0.2407008061506589
$ws.Range("P2").Value = 0.2407008061506589
$ws.Range("Q2").Value = 1.627429989916
$ws.Range("R2").Value = 14.646869909244
$ws.Range("S2").Value = 0.05381224180572483
$ws.Range("T2").Value = 0.05381224180572484
$ws.Range("G3").Value = 2.738607
$ws.Range("H3").Value = 8.215821
$ws.Range("I3").Value = 0.2235648590725649
$ws.Range("J3").Value = 0.223564859072565
$ws.Range("O3").Value = 0.4863657278767622
$ws.Range("P3").Value = 0.4863657278767622
$ws.Range("Q3").Value = 3.288423434355
$ws.Range("R3").Value = 29.595810909195
$ws.Range("S3").Value = 0.1087342854104938
$ws.Range("T3").Value = 0.1087342854104938
$ws.Range("G4").Value = 2.738607
$ws.Range("H4").Value = 8.215821
$ws.Range("I4").Value = 0.2235648590725649
$ws.Range("J4").Value = 0.223564859072565
$ws.Range("O4").Value = 0.2729334659725789
$ws.Range("P4").Value = 0.2729334659725789
$ws.Range("Q4").Value = 1.845361944893
$ws.Range("R4").Value = 16.608257504037
$ws.Range("S4").Value = 0.06101833185634629
$ws.Range("T4").Value = 0.0610183318563463
$ws.Range("I5").Value = 0.5119261090069511
$ws.Range("J5").Value = 0.5119261090069511
$ws.Range("M5").Value = 0.5942546666666666
$ws.Range("N5").Value = 1.782764
$ws.Range("O5").Value = 0.2407008061506589
$ws.Range("P5").Value = 0.2407008061506589
$ws.Range("Q5").Value = 3.726542292357777
$ws.Range("R5").Value = 33.53888063122
$ws.Range("S5").Value = 0.1232210271275432
$ws.Range("T5").Value = 0.1232210271275432
$ws.Range("I6").Value = 0.5119261090069511
$ws.Range("J6").Value = 0.5119261090069511
$ws.Range("O6").Value = 0.4863657278767622
$ws.Range("P6").Value = 0.4863657278767622
$ws.Range("S6").Value = 0.2489833146262845
$ws.Range("T6").Value = 0.2489833146262845
$ws.Range("I7").Value = 0.5119261090069511
$ws.Range("J7").Value = 0.5119261090069511
$ws.Range("O7").Value = 0.2729334659725789
$ws.Range("P7").Value = 0.2729334659725789
$ws.Range("S7").Value = 0.1397217672531234
$ws.Range("T7").Value = 0.1397217672531234
$ws.Range("I8").Value = 0.2645090319204839
$ws.Range("J8").Value = 0.2645090319204839
$ws.Range("M8").Value = 0.5942546666666666
$ws.Range("N8").Value = 1.782764
$ws.Range("O8").Value = 0.2407008061506589
$ws.Range("P8").Value = 0.2407008061506589
$ws.Range("Q8").Value = 1.925481191171111
$ws.Range("R8").Value = 17.32933072054
$ws.Range("S8").Value = 0.06366753721739084
$ws.Range("T8").Value = 0.06366753721739084
$ws.Range("I9").Value = 0.2645090319204839
$ws.Range("J9").Value = 0.2645090319204839
$ws.Range("O9").Value = 0.4863657278767622
$ws.Range("P9").Value = 0.4863657278767622
$ws.Range("S9").Value = 0.1286481278399839
$ws.Range("T9").Value = 0.1286481278399839
$ws.Range("I10").Value = 0.2645090319204839
$ws.Range("J10").Value = 0.2645090319204839
$ws.Range("O10").Value = 0.2729334659725789
$ws.Range("P10").Value = 0.2729334659725789
$ws.Range("S10").Value = 0.07219336686310916
$ws.Range("T10").Value = 0.07219336686310916
